$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.026.56"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").Value = "'2.963.76"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'595.48"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "'147.54"
$ws.Range("E6").Value = "  +1.51%  "
$ws.Range("D8").Value = "'2.962.70"
$ws.Range("E8").Value = "  +0.92%  "
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("E10").Value = "  +3.56%  "
$ws.Range("E11").Value = "  +6.20%  "
$ws.Range("D12").Value = "'0.446"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("E13").Value = "  +6.31%  "
$ws.Range("D14").Value = "'33.30"
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "'3.454.96"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").Value = "'62.895.33"
$ws.Range("E17").Value = "  +2.49%  "
$ws.Range("D18").Value = "'6.76"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "'2.974.39"
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("D20").Value = "'444.51"
$ws.Range("E20").Value = "  +2.36%  "
$ws.Range("D21").Value = "'13.51"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "'0.670"
$ws.Range("E22").Value = "  -1.69%  "
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").Value = "'11.29"
$ws.Range("E24").Value = "  +2.55%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'81.70"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").Value = "'2.15"
$ws.Range("E26").Value = "  -3.41%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'11.93"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").Value = "'7.36"
$ws.Range("E29").Value = "  +4.82%  "
$ws.Range("E31").Value = "  -3.32%  "
$ws.Range("D32").Value = "'0.0₃0975"
$ws.Range("E32").Value = "  +9.66%  "
$ws.Range("D33").Value = "'26.60"
$ws.Range("E33").Value = "  -0.71%  "
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'0.998"
$ws.Range("E36").Value = "  -1.58%  "
$ws.Range("D37").Value = "'3.14"
$ws.Range("E37").Value = "  +3.52%  "
$ws.Range("D38").Value = "'5.67"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").Value = "'2.07"
$ws.Range("E39").Value = "  +1.98%  "
$ws.Range("D40").Value = "'49.56"
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("D41").Value = "'8.55"
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("D42").Value = "'0.118"
$ws.Range("E42").Value = "  -4.82%  "
$ws.Range("D43").Value = "'0.284"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").Value = "'40.81"
$ws.Range("E44").Value = "  -4.42%  "
$ws.Range("D45").Value = "'2.717.72"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").Value = "'134.30"
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").Value = "'365.30"
$ws.Range("E47").Value = "  -1.94%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0339"
$ws.Range("E48").Value = "  -3.01%  "
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").Value = "'22.96"
$ws.Range("E51").Value = "  -4.54%  "
